# Notes to self.
#
# Insert two new todo bullets into the list:
#   1. "Unit tests for examples." as a new sub-bullet right before the
#      "Examples" bullet.
#   2. "Document minimum required instruction set and bump it in compiler
#      options (P4? Higher? What is the minimum for Vista or 7?)." as a new
#      top-level bullet right before the "Look for places where cleanup..."
#      bullet.
# Word's automatic "_GoBack" bookmark (marks the location of the most
# recent edit) is relocated so it sits at the start of the "Look for
# places where cleanup..." paragraph, matching where Word itself leaves it
# after the last edit made while authoring this change.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return 0
}

# --- Insertion 1: "Unit tests for examples." before "Examples" ---
$examplesIndex = Find-ParagraphIndex $d "Examples`r"
$examplesPara = $d.Paragraphs($examplesIndex)
$examplesPara.Range.InsertParagraphBefore()

$newPara1 = $d.Paragraphs($examplesIndex)
$newPara1.Range.ListFormat.ListLevelNumber = 2
$newPara1.Range.Text = "Unit tests for examples."

# --- Insertion 2: "Document minimum required instruction set..." before
#     "Look for places where cleanup..." ---
$cleanupIndex = Find-ParagraphIndex $d "Look for places where cleanup*"
$cleanupPara = $d.Paragraphs($cleanupIndex)
$cleanupPara.Range.InsertParagraphBefore()

$newPara2 = $d.Paragraphs($cleanupIndex)
$newPara2.Range.ListFormat.ListLevelNumber = 1
$newPara2.Range.Text = "Document minimum required instruction set and bump it in compiler options (P4? Higher? What is the minimum for Vista or 7?)."

# --- Relocate the "_GoBack" bookmark to mark the latest edit location,
#     i.e. the start of the (now shifted) "Look for places where
#     cleanup..." paragraph. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$cleanupIndex2 = Find-ParagraphIndex $d "Look for places where cleanup*"
$cleanupPara2 = $d.Paragraphs($cleanupIndex2)
$goBackRange = $cleanupPara2.Range.Duplicate
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
